$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A10").Value = 43744
$ws.Range("B10").Value = 240

$ws.Range("A11").Value = 43745
$ws.Range("B11").Value = 243

$ws.Range("B12").Select()
